$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: new entry for FOX / FOX924B (HCMOS output) ---
$ws.Range("A20").Value = "TCXO"
$ws.Range("B20").Value = "FOX"
$ws.Range("C20").Value = "FOX924B"
$ws.Range("D20").Value = 1.5
$ws.Range("E20").Value = 2500
$ws.Range("F20").Value = 300
$ws.Range("H20").ClearContents() | Out-Null
$ws.Range("I20").ClearContents() | Out-Null
$ws.Range("K20").Value = 6.41
$ws.Range("L20").Value = 3.3
$ws.Range("M20").Value = 6
$ws.Range("N20").ClearContents() | Out-Null
$ws.Range("O20").Value = "3.2x5"
$ws.Range("P20").Value = "HCMOS"
$ws.Range("J20").Value = "-30-85"

# --- Row 21: new entry for TXC Corp / 7N-20.000MBP-T (CMOS output) ---
$ws.Range("A21").Value = "TCXO"
$ws.Range("B21").Value = "TXC Corp"
$ws.Range("C21").Value = "7N-20.000MBP-T"
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 280
$ws.Range("F21").Value = 100
$ws.Range("H21").Value = 4.6
$ws.Range("I21").Value = 20
$ws.Range("J21").Value = "-40-85"
$ws.Range("K21").Value = 10.89
$ws.Range("L21").Value = 5
$ws.Range("M21").Value = 10
$ws.Range("N21").ClearContents() | Out-Null
$ws.Range("O21").Value = "7x5"
$ws.Range("P21").Value = "CMOS"

# --- Aging column (1ppm/yr) filled in last, across row 3, 20 and 21 ---
$ws.Range("G20").Value = "1ppm/yr"
$ws.Range("G21").Value = "1ppm/yr"
$ws.Range("G3").Value = "1ppm/yr"
$ws.Range("H3").ClearContents() | Out-Null
$ws.Range("I3").ClearContents() | Out-Null

# --- Apply row 20/21 cell formatting (matches the same pattern as row 14) ---
# Done after the values are set so number formats such as the quote-prefixed
# "Temp Range" column are not clobbered by the value assignment above.
$ws.Range("A14:P14").Copy() | Out-Null
$ws.Range("A20:P21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("H31").Select() | Out-Null
